# Refresh crypto price/volume data and reorder coin listing rows 7-18 and 46-47
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = '273.48'
$ws.Range("E2").Value = '-1.95%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = '26.67'
$ws.Range("E3").Value = '-2.28%'

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.24%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06313'
$ws.Range("E5").Value = '0.50%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = '6.871'
$ws.Range("E6").Value = '0.30%'

# Row 7
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = '1.200'
$ws.Range("E7").Value = '26.02%'

# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8703'
$ws.Range("E8").Value = '-0.67%'

# Row 9
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1458'
$ws.Range("E9").Value = '0.28%'

# Row 10
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05087'
$ws.Range("E10").Value = '-2.23%'

# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07377'
$ws.Range("E11").Value = '1.51%'

# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03032'
$ws.Range("E12").Value = '-3.24%'

# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09043'
$ws.Range("E13").Value = '-0.06%'

# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001576'
$ws.Range("E14").Value = '1.04%'

# Row 15
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006325'
$ws.Range("E15").Value = '0.90%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006009'
$ws.Range("E16").Value = '-0.12%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value = '3.454'
$ws.Range("E17").Value = '0.05%'

# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = '3.346'
$ws.Range("E18").Value = '2.26%'

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.21%'

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.54%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1325'
$ws.Range("E21").Value = '1.18%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = '3.922'
$ws.Range("E22").Value = '2.15%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04406'
$ws.Range("E23").Value = '2.13%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001176'
$ws.Range("E24").Value = '0.39%'

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.21%'

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.52%'

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.15%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006718'
$ws.Range("E41").Value = '0.24%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1164'
$ws.Range("E42").Value = '1.06%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002100'
$ws.Range("E43").Value = '0.07%'

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01253'
$ws.Range("E44").Value = '-9.26%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005311'
$ws.Range("E45").Value = '2.83%'

# Row 46
$ws.Range("B46").Value = 'CoinbaseStockToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D46").Value = '0.02000'
$ws.Range("E46").Value = '-33.03%'

# Row 47
$ws.Range("B47").Value = 'BOLO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D47").Value = '2.988'
$ws.Range("E47").Value = '28.26%'
